$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 16.820675
$ws.Range("H2").Value = 50.462025
$ws.Range("I2").Value = 0.8427583848046372
$ws.Range("J2").Value = 0.8427583848046373
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 87.94127933333334
$ws.Range("N2").Value = 263.823838
$ws.Range("O2").Value = 0.4109331243514438
$ws.Range("P2").Value = 0.4109331243514437
$ws.Range("Q2").Value = 1479.231678750217
$ws.Range("R2").Value = 13313.08510875195
$ws.Range("S2").Value = 0.3463173361411459
$ws.Range("T2").Value = 0.3463173361411459
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 16.820675
$ws.Range("H3").Value = 50.462025
$ws.Range("I3").Value = 0.8427583848046372
$ws.Range("J3").Value = 0.8427583848046373
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 52.441971
$ws.Range("N3").Value = 157.325913
$ws.Range("O3").Value = 0.2450515065683088
$ws.Range("P3").Value = 0.2450515065683087
$ws.Range("Q3").Value = 882.1093505504249
$ws.Range("R3").Value = 7938.984154953825
$ws.Range("S3").Value = 0.2065192118694508
$ws.Range("T3").Value = 0.2065192118694508
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 16.820675
$ws.Range("H4").Value = 50.462025
$ws.Range("I4").Value = 0.8427583848046372
$ws.Range("J4").Value = 0.8427583848046373
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 54.667459
$ws.Range("N4").Value = 164.002377
$ws.Range("O4").Value = 0.255450795093328
$ws.Range("P4").Value = 0.255450795093328
$ws.Range("Q4").Value = 919.5435609148249
$ws.Range("R4").Value = 8275.892048233425
$ws.Range("S4").Value = 0.2152832994699135
$ws.Range("T4").Value = 0.2152832994699134
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 16.820675
$ws.Range("H5").Value = 50.462025
$ws.Range("I5").Value = 0.8427583848046372
$ws.Range("J5").Value = 0.8427583848046373
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 18.95316166666667
$ws.Range("N5").Value = 56.85948500000001
$ws.Range("O5").Value = 0.08856457398691947
$ws.Range("P5").Value = 0.08856457398691944
$ws.Range("Q5").Value = 318.8049726174584
$ws.Range("R5").Value = 2869.244753557125
$ws.Range("S5").Value = 0.07463853732412704
$ws.Range("T5").Value = 0.07463853732412702
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.464483
$ws.Range("H6").Value = 4.393449
$ws.Range("I6").Value = 0.07337430439942808
$ws.Range("J6").Value = 0.07337430439942808
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 87.94127933333334
$ws.Range("N6").Value = 263.823838
$ws.Range("O6").Value = 0.4109331243514438
$ws.Range("P6").Value = 0.4109331243514437
$ws.Range("Q6").Value = 128.788508581918
$ws.Range("R6").Value = 1159.096577237262
$ws.Range("S6").Value = 0.03015193215397087
$ws.Range("T6").Value = 0.03015193215397086
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.464483
$ws.Range("H7").Value = 4.393449
$ws.Range("I7").Value = 0.07337430439942808
$ws.Range("J7").Value = 0.07337430439942808
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 52.441971
$ws.Range("N7").Value = 157.325913
$ws.Range("O7").Value = 0.2450515065683088
$ws.Range("P7").Value = 0.2450515065683087
$ws.Range("Q7").Value = 76.80037501599301
$ws.Range("R7").Value = 691.2033751439371
$ws.Range("S7").Value = 0.01798048383648154
$ws.Range("T7").Value = 0.01798048383648153
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.464483
$ws.Range("H8").Value = 4.393449
$ws.Range("I8").Value = 0.07337430439942808
$ws.Range("J8").Value = 0.07337430439942808
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 54.667459
$ws.Range("N8").Value = 164.002377
$ws.Range("O8").Value = 0.255450795093328
$ws.Range("P8").Value = 0.255450795093328
$ws.Range("Q8").Value = 80.05956435869702
$ws.Range("R8").Value = 720.536079228273
$ws.Range("S8").Value = 0.01874352439825378
$ws.Range("T8").Value = 0.01874352439825377
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.464483
$ws.Range("H9").Value = 4.393449
$ws.Range("I9").Value = 0.07337430439942808
$ws.Range("J9").Value = 0.07337430439942808
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 18.95316166666667
$ws.Range("N9").Value = 56.85948500000001
$ws.Range("O9").Value = 0.08856457398691947
$ws.Range("P9").Value = 0.08856457398691944
$ws.Range("Q9").Value = 27.75658305708501
$ws.Range("R9").Value = 249.8092475137651
$ws.Range("S9").Value = 0.006498364010721899
$ws.Range("T9").Value = 0.006498364010721897
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.1122456666666666
$ws.Range("H10").Value = 0.336737
$ws.Range("I10").Value = 0.005623791954919746
$ws.Range("J10").Value = 0.005623791954919746
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 87.94127933333334
$ws.Range("N10").Value = 263.823838
$ws.Range("O10").Value = 0.4109331243514438
$ws.Range("P10").Value = 0.4109331243514437
$ws.Range("Q10").Value = 9.871027526289554
$ws.Range("R10").Value = 88.839247736606
$ws.Range("S10").Value = 0.002311002398737685
$ws.Range("T10").Value = 0.002311002398737685
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.1122456666666666
$ws.Range("H11").Value = 0.336737
$ws.Range("I11").Value = 0.005623791954919746
$ws.Range("J11").Value = 0.005623791954919746
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 52.441971
$ws.Range("N11").Value = 157.325913
$ws.Range("O11").Value = 0.2450515065683088
$ws.Range("P11").Value = 0.2450515065683087
$ws.Range("Q11").Value = 5.886383996208999
$ws.Range("R11").Value = 52.977455965881
$ws.Range("S11").Value = 0.001378118691179818
$ws.Range("T11").Value = 0.001378118691179818
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.1122456666666666
$ws.Range("H12").Value = 0.336737
$ws.Range("I12").Value = 0.005623791954919746
$ws.Range("J12").Value = 0.005623791954919746
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 54.667459
$ws.Range("N12").Value = 164.002377
$ws.Range("O12").Value = 0.255450795093328
$ws.Range("P12").Value = 0.255450795093328
$ws.Range("Q12").Value = 6.136185380427666
$ws.Range("R12").Value = 55.22566842384899
$ws.Range("S12").Value = 0.001436602126323711
$ws.Range("T12").Value = 0.00143660212632371
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.1122456666666666
$ws.Range("H13").Value = 0.336737
$ws.Range("I13").Value = 0.005623791954919746
$ws.Range("J13").Value = 0.005623791954919746
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 18.95316166666667
$ws.Range("N13").Value = 56.85948500000001
$ws.Range("O13").Value = 0.08856457398691947
$ws.Range("P13").Value = 0.08856457398691944
$ws.Range("Q13").Value = 2.127410266716111
$ws.Range("R13").Value = 19.146692400445
$ws.Range("S13").Value = 0.0004980687386785324
$ws.Range("T13").Value = 0.0004980687386785321
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 1.561668
$ws.Range("H14").Value = 4.685003999999999
$ws.Range("I14").Value = 0.07824351884101489
$ws.Range("J14").Value = 0.07824351884101489
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 87.94127933333334
$ws.Range("N14").Value = 263.823838
$ws.Range("O14").Value = 0.4109331243514438
$ws.Range("P14").Value = 0.4109331243514437
$ws.Range("Q14").Value = 137.335081813928
$ws.Range("R14").Value = 1236.015736325352
$ws.Range("S14").Value = 0.03215285365758931
$ws.Range("T14").Value = 0.0321528536575893
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 1.561668
$ws.Range("H15").Value = 4.685003999999999
$ws.Range("I15").Value = 0.07824351884101489
$ws.Range("J15").Value = 0.07824351884101489
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 52.441971
$ws.Range("N15").Value = 157.325913
$ws.Range("O15").Value = 0.2450515065683088
$ws.Range("P15").Value = 0.2450515065683087
$ws.Range("Q15").Value = 81.89694796762799
$ws.Range("R15").Value = 737.0725317086519
$ws.Range("S15").Value = 0.01917369217119655
$ws.Range("T15").Value = 0.01917369217119655
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 1.561668
$ws.Range("H16").Value = 4.685003999999999
$ws.Range("I16").Value = 0.07824351884101489
$ws.Range("J16").Value = 0.07824351884101489
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 54.667459
$ws.Range("N16").Value = 164.002377
$ws.Range("O16").Value = 0.255450795093328
$ws.Range("P16").Value = 0.255450795093328
$ws.Range("Q16").Value = 85.372421361612
$ws.Range("R16").Value = 768.3517922545078
$ws.Range("S16").Value = 0.01998736909883705
$ws.Range("T16").Value = 0.01998736909883704
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 1.561668
$ws.Range("H17").Value = 4.685003999999999
$ws.Range("I17").Value = 0.07824351884101489
$ws.Range("J17").Value = 0.07824351884101489
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 18.95316166666667
$ws.Range("N17").Value = 56.85948500000001
$ws.Range("O17").Value = 0.08856457398691947
$ws.Range("P17").Value = 0.08856457398691944
$ws.Range("Q17").Value = 29.59854607366
$ws.Range("R17").Value = 266.38691466294
$ws.Range("S17").Value = 0.006929603913391991
$ws.Range("T17").Value = 0.006929603913391989
